$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''51.513.29'
$ws.Range("E2").Value = '  -0.16%  '

$ws.Range("D3").Value = '''3.100.11'
$ws.Range("E3").Value = '  +2.49%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '''384.24'
$ws.Range("E5").Value = '  +1.34%  '

$ws.Range("D6").Value = '''103.21'
$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("D7").Value = '''0.540'
$ws.Range("E7").Value = '  -0.82%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = '''0.585'
$ws.Range("E9").Value = '  -1.83%  '

$ws.Range("D10").Value = '''36.94'
$ws.Range("E10").Value = '  +0.48%  '

$ws.Range("E11").Value = '  -0.08%  '

$ws.Range("D12").Value = '''0.0856'
$ws.Range("E12").Value = '  -0.51%  '

$ws.Range("D13").Value = '''3.587.39'
$ws.Range("E13").Value = '  +2.54%  '

$ws.Range("D14").Value = '''18.64'
$ws.Range("E14").Value = '  +0.64%  '

$ws.Range("D15").Value = '''7.84'
$ws.Range("E15").Value = '  +1.06%  '

$ws.Range("D16").Value = '''3.100.01'
$ws.Range("E16").Value = '  +2.14%  '

$ws.Range("D17").Value = '''11.10'
$ws.Range("E17").Value = '  +5.97%  '

$ws.Range("D18").Value = '''0.996'
$ws.Range("E18").Value = '  +1.28%  '

$ws.Range("D19").Value = '''51.556.62'
$ws.Range("E19").Value = '  -0.08%  '

$ws.Range("D20").Value = '''3.34'
$ws.Range("E20").Value = '  +9.57%  '

$ws.Range("D21").Value = '''12.39'
$ws.Range("E21").Value = '  -0.71%  '

$ws.Range("E22").Value = '  +0.08%  '

$ws.Range("D23").Value = '''70.00'
$ws.Range("E23").Value = '  -0.09%  '

$ws.Range("D24").Value = '''266.25'
$ws.Range("E24").Value = '  -0.70%  '

$ws.Range("E25").Value = '  +0.01%  '

$ws.Range("D26").Value = '''8.15'
$ws.Range("E26").Value = '  -0.54%  '

$ws.Range("D27").Value = '''27.07'
$ws.Range("E27").Value = '  +3.25%  '

$ws.Range("D28").Value = '''7.26'
$ws.Range("E28").Value = '  -3.52%  '

$ws.Range("E29").Value = '  +0.07%  '

$ws.Range("D30").Value = '''0.168'
$ws.Range("E30").Value = '  -3.04%  '

$ws.Range("E31").Value = '  -2.11%  '

$ws.Range("D32").Value = '''10.34'
$ws.Range("E32").Value = '  +0.18%  '

$ws.Range("D33").Value = '''35.42'
$ws.Range("E33").Value = '  +3.14%  '

$ws.Range("D34").Value = '''0.0467'
$ws.Range("E34").Value = '  +2.93%  '

$ws.Range("B35").Value = 'Toncoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D35").Value = '''2.07'
$ws.Range("E35").Value = '  +0.60%  '

$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D36").Value = '''50.29'
$ws.Range("E36").Value = '  -0.49%  '

$ws.Range("D37").Value = '''0.999'
$ws.Range("E37").Value = '  -0.11%  '

$ws.Range("E38").Value = '  +2.34%  '

$ws.Range("D39").Value = '''0.301'
$ws.Range("E39").Value = '  +4.32%  '

$ws.Range("E40").Value = '  +1.01%  '

$ws.Range("D41").Value = '''128.92'
$ws.Range("E41").Value = '  +1.67%  '

$ws.Range("D42").Value = '''16.58'
$ws.Range("E42").Value = '  -4.91%  '

$ws.Range("E43").Value = '  -0.48%  '

$ws.Range("E44").Value = '  -2.34%  '

$ws.Range("D45").Value = '''22.41'
$ws.Range("E45").Value = '  +2.14%  '

$ws.Range("D46").Value = '''3.66'
$ws.Range("E46").Value = '  -1.72%  '

$ws.Range("E47").Value = '  +3.83%  '

$ws.Range("D48").Value = '''2.08'
$ws.Range("E48").Value = '  +1.25%  '

$ws.Range("D49").Value = '''2.060.38'
$ws.Range("E49").Value = '  +1.30%  '

$ws.Range("D50").Value = '''0.0332'
$ws.Range("E50").Value = '  +3.47%  '

$ws.Range("E51").Value = '  +12.86%  '
